$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the chart object and the (pre-insert) default width of the
# column we are about to insert, so we can shift the chart anchor by
# exactly one column afterwards (this sandboxed engine does not re-anchor
# floating objects automatically on a column insert).
$co = $ws.ChartObjects().Item(1)
$newColDefaultWidth = $ws.Columns("G").Width

# Insert a new column before G - this pushes the old "Places/Encounters"
# summary block (H:J) one column to the right (I:K) and lets Excel's
# formula engine adjust all the relative formulas automatically.
$ws.Columns("G").Insert()

# Re-anchor the chart so it keeps sitting to the right of the summary
# block, exactly one column further right than before.
$co.Left = $co.Left + $newColDefaultWidth

# Give the new "story" column a width (the old best-fit summary-label
# column slides from H to I on its own, keeping its existing width).
$ws.Columns("G").ColumnWidth = 12.375

# Fill in the new "story" labels next to each encounter row (1..11), in
# the same order they were authored in (matches the shared-string table
# order of the original edit).
$ws.Range("G7").Value = "Astrakan"
$ws.Range("G8").Value = "Krasnovodsk"
$ws.Range("G9").Value = "Turcomans passing by"
$ws.Range("G10").Value = "Meeting a thirsty turcoman"
$ws.Range("G15").Value = "Turcoman selling supplies"
$ws.Range("G13").Value = "Turcoman bandit raid"
$ws.Range("G14").Value = "Extreme heat"
$ws.Range("G17").Value = "Major Frankenburg verge of death"
$ws.Range("G16").Value = "Oasis Mirage"
$ws.Range("G12").Value = "Real Oasis"
$ws.Range("G11").Value = "Strong heat causing thirst"

# Re-order a handful of the Sink(E)/Gain(F) values to line up with the
# newly organised story order.
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 3
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 6
$ws.Range("E13").Value = 7
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 0
